$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.74%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.45%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.222"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.48%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07256"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "8.24%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.796"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.12%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.753"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "9.05%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.457"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.58%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9060"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01643"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2,444.35%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1681"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.95%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07402"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "8.93%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07937"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.80%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02968"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.30%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09906"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "10.29%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001498"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-6.16%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04559"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.84%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006453"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.54%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.481"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.93%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.227"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.06%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3334"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.24%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1328"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.31%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.276"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.01%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1636"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.49%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001227"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.49%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004413"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "7.00%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001307"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "9.11%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "8.34%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04464"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.55%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007057"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.02%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1337"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.79%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002353"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.21%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01278"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.29%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006107"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.58%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.83%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01616"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "7.38%"
